# Lesson 16 updates: clear the leftover "#NULL!" placeholder cells from the
# bottom of the singer-heights table (rows 22-40) and reset the sheet
# selection back to the default (A1) instead of the stale D1 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D22:D36 only had a stray "#NULL!" text value in column D - clear those.
$ws.Range("D22:D36").ClearContents()

# Row 37 only has real data in A:B; C37/D37 held "#NULL!" placeholders.
$ws.Range("C37:D37").ClearContents()

# Rows 38-40 only have real data in column A; B:D held "#NULL!" placeholders.
$ws.Range("B38:D40").ClearContents()

# Restore the default selection (A1) on the active sheet instead of D1.
$ws.Range("A1").Select()
